$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new trip rows (14-30), mirroring the existing data pattern:
# column A = departure_location, B = arrival_location, C = departure_date,
# D = arrival_date (only present on a handful of rows, same as source data).
$ws.Cells.Item(14, 1).Value = "São Paulo"
$ws.Cells.Item(14, 2).Value = "Rio de Janeiro"
$ws.Cells.Item(14, 3).Value = 45079
$ws.Cells.Item(14, 3).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(15, 1).Value = "São Paulo"
$ws.Cells.Item(15, 2).Value = "Rio de Janeiro"
$ws.Cells.Item(15, 3).Value = 45079
$ws.Cells.Item(15, 3).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(16, 1).Value = "São Paulo"
$ws.Cells.Item(16, 2).Value = "Rio de Janeiro"
$ws.Cells.Item(16, 3).Value = 45080
$ws.Cells.Item(16, 3).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(17, 1).Value = "São Paulo"
$ws.Cells.Item(17, 2).Value = "Rio de Janeiro"
$ws.Cells.Item(17, 3).Value = 45080
$ws.Cells.Item(17, 3).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(18, 1).Value = "Florianópolis"
$ws.Cells.Item(18, 2).Value = "Belo Horizonte"
$ws.Cells.Item(18, 3).Value = 45081
$ws.Cells.Item(18, 3).NumberFormat = "dd/mm/yy"
$ws.Cells.Item(18, 4).Value = 45091
$ws.Cells.Item(18, 4).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(19, 1).Value = "Florianópolis"
$ws.Cells.Item(19, 2).Value = "Rio de Janeiro"
$ws.Cells.Item(19, 3).Value = 45081
$ws.Cells.Item(19, 3).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(20, 1).Value = "Florianópolis"
$ws.Cells.Item(20, 2).Value = "São Paulo"
$ws.Cells.Item(20, 3).Value = 45082
$ws.Cells.Item(20, 3).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(21, 1).Value = "Florianópolis"
$ws.Cells.Item(21, 2).Value = "Curitiba"
$ws.Cells.Item(21, 3).Value = 45083
$ws.Cells.Item(21, 3).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(22, 1).Value = "Rio de Janeiro"
$ws.Cells.Item(22, 2).Value = "São Paulo"
$ws.Cells.Item(22, 3).Value = 45084
$ws.Cells.Item(22, 3).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(23, 1).Value = "Rio de Janeiro"
$ws.Cells.Item(23, 2).Value = "Recife"
$ws.Cells.Item(23, 3).Value = 45085
$ws.Cells.Item(23, 3).NumberFormat = "dd/mm/yy"
$ws.Cells.Item(23, 4).Value = 45091
$ws.Cells.Item(23, 4).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(24, 1).Value = "Rio de Janeiro"
$ws.Cells.Item(24, 2).Value = "Recife"
$ws.Cells.Item(24, 3).Value = 45086
$ws.Cells.Item(24, 3).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(25, 1).Value = "Rio de Janeiro"
$ws.Cells.Item(25, 2).Value = "Belo Horizonte"
$ws.Cells.Item(25, 3).Value = 45087
$ws.Cells.Item(25, 3).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(26, 1).Value = "Belo Horizonte"
$ws.Cells.Item(26, 2).Value = "Recife"
$ws.Cells.Item(26, 3).Value = 45088
$ws.Cells.Item(26, 3).NumberFormat = "dd/mm/yy"
$ws.Cells.Item(26, 4).Value = 45091
$ws.Cells.Item(26, 4).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(27, 1).Value = "Belo Horizonte"
$ws.Cells.Item(27, 2).Value = "Rio de Janeiro"
$ws.Cells.Item(27, 3).Value = 45089
$ws.Cells.Item(27, 3).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(28, 1).Value = "Belo Horizonte"
$ws.Cells.Item(28, 2).Value = "São Paulo"
$ws.Cells.Item(28, 3).Value = 45090
$ws.Cells.Item(28, 3).NumberFormat = "dd/mm/yy"
$ws.Cells.Item(28, 4).Value = 45091
$ws.Cells.Item(28, 4).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(29, 1).Value = "Belo Horizonte"
$ws.Cells.Item(29, 2).Value = "Curitiba"
$ws.Cells.Item(29, 3).Value = 45091
$ws.Cells.Item(29, 3).NumberFormat = "dd/mm/yy"

$ws.Cells.Item(30, 1).Value = "Belo Horizonte"
$ws.Cells.Item(30, 2).Value = "São Paulo"
$ws.Cells.Item(30, 3).Value = 45092
$ws.Cells.Item(30, 3).NumberFormat = "dd/mm/yy"

# Restore the active-cell selection to where editing left off.
$ws.Range("I7").Select()